# Update the US03 (row 4) and US04 (row 5) detail rows on the "Sprint1" sheet
# with the DateValidation.py / TestDateValidation.py source & test function
# references, per the commit "update the US03 and US04's detail in TeamDReport".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# xlLeft
$xlLeft = -4131

# ---- Row 4 (US03 - Birth before death) ----
$ws.Range("G4").Value = 15
$ws.Range("I4").Value = 42047

$ws.Range("K4").Value = "DateValidation.py"

$ws.Range("L4").Value = "validate_birth_before_death"
$ws.Range("L4").WrapText = $true
$ws.Range("L4").HorizontalAlignment = $xlLeft

$ws.Range("M4").Value = 15

$ws.Range("O4").Value = "TestDateValidation.py"

$ws.Range("P4").Value = "test_validate_birth_before_death_XXX"
$ws.Range("P4").WrapText = $true
$ws.Range("P4").HorizontalAlignment = $xlLeft

$ws.Range("Q4").Value = 21

$ws.Rows.Item(4).RowHeight = 28

# ---- Row 5 (US04 - Marriage before divorce) ----
$ws.Range("G5").Value = 14
$ws.Range("I5").Value = 42047

$ws.Range("K5").Value = "DateValidation.py"

$ws.Range("L5").Value = "validate_marraige_before_divorce"
$ws.Range("L5").WrapText = $true
$ws.Range("L5").HorizontalAlignment = $xlLeft

$ws.Range("M5").Value = 14

$ws.Range("O5").Value = "TestDateValidation.py"

$ws.Range("P5").Value = "test_validate_marraige_before_divorce_XXX"
$ws.Range("P5").WrapText = $true
$ws.Range("P5").HorizontalAlignment = $xlLeft

$ws.Range("Q5").Value = 20

$ws.Rows.Item(5).RowHeight = 28

# ---- Restore sheet view selection (matches author's last cursor position) ----
$ws.Range("O12").Select() | Out-Null
